$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 22
$ws.Range("H2").Value = 0.005909369219200378
$ws.Range("I2").Value = 0.005909369219200378
$ws.Range("J2").Value = 0.4879686802983292
$ws.Range("K2").Value = 0.4879686802983292
$ws.Range("L2").Value = 36.18742926255743
$ws.Range("M2").Value = "[11.380987989867677, 60.99387053524718]"
$ws.Range("N2").Value = 0.005191697364629544
$ws.Range("O2").Value = 0.005191697364629544
$ws.Range("P2").Value = 1.515763422452733
$ws.Range("Q2").Value = "[0.5849211547224238, 2.446605690183042]"
$ws.Range("R2").Value = 0.002010056663019233
$ws.Range("S2").Value = 0.002010056663019233
$ws.Range("T2").Value = 56.27270192096265
$ws.Range("U2").Value = "[41.59607485235995, 70.94932898956536]"
$ws.Range("V2").Value = [double]"8.76527295190499e-10"
$ws.Range("W2").Value = [double]"8.76527295190499e-10"
$ws.Range("X2").Value = 16.69269269269269
$ws.Range("Y2").Value = 13.43343343343343
$ws.Range("Z2").Value = 19.95195195195195

# Row 3
$ws.Range("F3").Value = 22
$ws.Range("H3").Value = 0.001888399901640248
$ws.Range("I3").Value = 0.001888399901640248
$ws.Range("J3").Value = 0.0007030343140426165
$ws.Range("K3").Value = 0.0007030343140426165
$ws.Range("L3").Value = 38.86707355554678
$ws.Range("M3").Value = "[12.675312544068916, 65.05883456702465]"
$ws.Range("N3").Value = 0.004525458792056902
$ws.Range("O3").Value = 0.004525458792056902
$ws.Range("P3").Value = 2.056658253701426
$ws.Range("Q3").Value = "[1.3270791789938858, 2.7862373284089657]"
$ws.Range("R3").Value = [double]"9.373189608918864e-07"
$ws.Range("S3").Value = [double]"9.373189608918864e-07"
$ws.Range("T3").Value = 47.87919110330442
$ws.Range("U3").Value = "[33.37982171385194, 62.378560492756904]"
$ws.Range("V3").Value = [double]"3.362543821161523e-08"
$ws.Range("W3").Value = [double]"3.362543821161523e-08"
$ws.Range("X3").Value = 14.7987987987988
$ws.Range("Y3").Value = 12.24424424424424
$ws.Range("Z3").Value = 17.35335335335336

# Row 4
$ws.Range("F4").Value = 22
$ws.Range("H4").Value = 0.001912418536416105
$ws.Range("I4").Value = 0.001912418536416105
$ws.Range("J4").Value = 0.009680995536788806
$ws.Range("K4").Value = 0.009680995536788806
$ws.Range("L4").Value = 37.50472895389199
$ws.Range("M4").Value = "[11.162759720584368, 63.846698187199614]"
$ws.Range("N4").Value = 0.006273080906055384
$ws.Range("O4").Value = 0.006273080906055384
$ws.Range("P4").Value = 2.572395185822273
$ws.Range("Q4").Value = "[1.8805529598065034, 3.2642374118380424]"
$ws.Range("R4").Value = [double]"1.931151460965452e-09"
$ws.Range("S4").Value = [double]"1.931151460965452e-09"
$ws.Range("T4").Value = 56.1261500574043
$ws.Range("U4").Value = "[41.655275576280374, 70.59702453852822]"
$ws.Range("V4").Value = [double]"6.485363357455753e-10"
$ws.Range("W4").Value = [double]"6.485363357455753e-10"
$ws.Range("X4").Value = 12.99299299299299
$ws.Range("Y4").Value = 10.57057057057057
$ws.Range("Z4").Value = 15.41541541541541

# Row 5
$ws.Range("F5").Value = 22
$ws.Range("H5").Value = 0.0003410248289880125
$ws.Range("I5").Value = 0.0003410248289880125
$ws.Range("J5").Value = 0.5845104513583419
$ws.Range("K5").Value = 0.5845104513583419
$ws.Range("L5").Value = 42.62279883054553
$ws.Range("M5").Value = "[19.27329027017754, 65.97230739091353]"
$ws.Range("N5").Value = 0.000627685800024036
$ws.Range("O5").Value = 0.000627685800024036
$ws.Range("P5").Value = -3.018947895341543
$ws.Range("Q5").Value = "[-3.723369070921236, -2.31452671976185]"
$ws.Range("R5").Value = [double]"4.225197969276451e-11"
$ws.Range("S5").Value = [double]"4.225197969276451e-11"
$ws.Range("T5").Value = 58.42106281062149
$ws.Range("U5").Value = "[44.134743654201415, 72.70738196704157]"
$ws.Range("V5").Value = [double]"1.565549911930475e-10"
$ws.Range("W5").Value = [double]"1.565549911930475e-10"
$ws.Range("X5").Value = 10.57057057057057
$ws.Range("Y5").Value = 8.104104104104106
$ws.Range("Z5").Value = 13.03703703703703

# Row 6
$ws.Range("F6").Value = 24.69000000000042
$ws.Range("H6").Value = [double]"3.982212707243082e-05"
$ws.Range("I6").Value = [double]"3.982212707243082e-05"
$ws.Range("L6").Value = 52.73830807943239
$ws.Range("M6").Value = "[26.86213817890777, 78.614477979957]"
$ws.Range("N6").Value = 0.0001680132856138883
$ws.Range("O6").Value = 0.0001680132856138883
$ws.Range("P6").Value = 2.899447874484274
$ws.Range("Q6").Value = "[2.320816194543811, 3.4780795544247365]"
$ws.Range("R6").Value = [double]"3.894662370385049e-13"
$ws.Range("S6").Value = [double]"3.894662370385049e-13"
$ws.Range("T6").Value = 71.31001424970292
$ws.Range("U6").Value = "[56.62143089384324, 85.9985976055626]"
$ws.Range("V6").Value = [double]"1.045608044591972e-12"
$ws.Range("W6").Value = [double]"1.045608044591972e-12"
$ws.Range("X6").Value = 13.29651651651674
$ws.Range("Y6").Value = 11.02276276276295
$ws.Range("Z6").Value = 15.57027027027054

# Row 7
$ws.Range("F7").Value = 24.69000000000042
$ws.Range("H7").Value = 0.0001319474139096499
$ws.Range("I7").Value = 0.0001319474139096499
$ws.Range("L7").Value = 52.52700522177441
$ws.Range("M7").Value = "[24.86500337313943, 80.18900707040939]"
$ws.Range("N7").Value = 0.0004007714991238753
$ws.Range("O7").Value = 0.0004007714991238753
$ws.Range("P7").Value = 2.974921571867812
$ws.Range("Q7").Value = "[2.371131992799503, 3.5787111509361207]"
$ws.Range("R7").Value = [double]"6.608047442568932e-13"
$ws.Range("S7").Value = [double]"6.608047442568932e-13"
$ws.Range("T7").Value = 59.50530416724568
$ws.Range("U7").Value = "[43.7244227401106, 75.28618559438075]"
$ws.Range("V7").Value = [double]"1.349655720517262e-09"
$ws.Range("W7").Value = [double]"1.349655720517262e-09"
$ws.Range("X7").Value = 12.99993993994016
$ws.Range("Y7").Value = 10.62732732732751
$ws.Range("Z7").Value = 15.37255255255282

# Row 8
$ws.Range("F8").Value = 24.69000000000042
$ws.Range("H8").Value = 0.003328840570210101
$ws.Range("I8").Value = 0.003328840570210101
$ws.Range("L8").Value = 42.30006773236263
$ws.Range("M8").Value = "[13.234835407829038, 71.36530005689623]"
$ws.Range("N8").Value = 0.005289701498068622
$ws.Range("O8").Value = 0.005289701498068622
$ws.Range("P8").Value = 2.912026824048196
$ws.Range("Q8").Value = "[2.094395102393195, 3.729658545703198]"
$ws.Range("R8").Value = [double]"5.642311950992962e-09"
$ws.Range("S8").Value = [double]"5.642311950992962e-09"
$ws.Range("T8").Value = 58.37279476423916
$ws.Range("U8").Value = "[41.82785391741224, 74.91773561106608]"
$ws.Range("V8").Value = [double]"7.096496057457102e-09"
$ws.Range("W8").Value = [double]"7.096496057457102e-09"
$ws.Range("X8").Value = 13.24708708708731
$ws.Range("Y8").Value = 10.03417417417434
$ws.Range("Z8").Value = 16.46000000000028

# Row 9
$ws.Range("F9").Value = 24.69000000000042
$ws.Range("H9").Value = [double]"8.621976173572854e-06"
$ws.Range("I9").Value = [double]"8.621976173572854e-06"
$ws.Range("L9").Value = 57.27743595286707
$ws.Range("M9").Value = "[28.79859785170143, 85.7562740540327]"
$ws.Range("N9").Value = 0.000199076144605792
$ws.Range("O9").Value = 0.000199076144605792
$ws.Range("P9").Value = 2.333395144107734
$ws.Range("Q9").Value = "[1.8553950606786556, 2.8113952275368126]"
$ws.Range("R9").Value = [double]"8.817391261572993e-13"
$ws.Range("S9").Value = [double]"8.817391261572993e-13"
$ws.Range("T9").Value = 75.45295088164715
$ws.Range("U9").Value = "[60.61737816796055, 90.28852359533376]"
$ws.Range("V9").Value = [double]"2.433608869978343e-13"
$ws.Range("W9").Value = [double]"2.433608869978343e-13"
$ws.Range("X9").Value = 15.52084084084111
$ws.Range("Y9").Value = 13.64252252252275
$ws.Range("Z9").Value = 17.39915915915946

# Row 10
$ws.Range("F10").Value = 24.69000000000042
$ws.Range("H10").Value = [double]"2.278237908248659e-05"
$ws.Range("I10").Value = [double]"2.278237908248659e-05"
$ws.Range("L10").Value = 60.48600577093592
$ws.Range("M10").Value = "[30.244837739410315, 90.72717380246152]"
$ws.Range("N10").Value = 0.0002134816420458208
$ws.Range("O10").Value = 0.0002134816420458208
$ws.Range("P10").Value = 2.119553001521041
$ws.Range("Q10").Value = "[1.6038160694001933, 2.6352899336418885]"
$ws.Range("R10").Value = [double]"1.364968138517497e-10"
$ws.Range("S10").Value = [double]"1.364968138517497e-10"
$ws.Range("T10").Value = 63.16950542219755
$ws.Range("U10").Value = "[46.545555347638434, 79.79345549675665]"
$ws.Range("V10").Value = [double]"1.10643516570974e-09"
$ws.Range("W10").Value = [double]"1.10643516570974e-09"
$ws.Range("X10").Value = 16.36114114114142
$ws.Range("Y10").Value = 14.33453453453478
$ws.Range("Z10").Value = 18.38774774774807

# Row 11
$ws.Range("F11").Value = 24.69000000000042
$ws.Range("H11").Value = 0.0007955063531688289
$ws.Range("I11").Value = 0.0007955063531688289
$ws.Range("L11").Value = 53.13673858752052
$ws.Range("M11").Value = "[16.94355972694818, 89.32991744809286]"
$ws.Range("N11").Value = 0.00493396630749543
$ws.Range("O11").Value = 0.00493396630749543
$ws.Range("P11").Value = 2.207605648468503
$ws.Range("Q11").Value = "[1.5786581702723481, 2.8365531266646578]"
$ws.Range("R11").Value = [double]"8.03804067750491e-09"
$ws.Range("S11").Value = [double]"8.03804067750491e-09"
$ws.Range("T11").Value = 72.22239529531593
$ws.Range("U11").Value = "[53.598480825769684, 90.84630976486218]"
$ws.Range("V11").Value = [double]"6.512732575458813e-10"
$ws.Range("W11").Value = [double]"6.512732575458813e-10"
$ws.Range("X11").Value = 16.01513513513541
$ws.Range("Y11").Value = 13.5436636636639
$ws.Range("Z11").Value = 18.48660660660692

# Row 12
$ws.Range("F12").Value = 24.69000000000042
$ws.Range("H12").Value = [double]"5.738323242332477e-05"
$ws.Range("I12").Value = [double]"5.738323242332477e-05"
$ws.Range("L12").Value = 51.62411878688288
$ws.Range("M12").Value = "[22.58158684765587, 80.66665072610989]"
$ws.Range("N12").Value = 0.0008374991573700807
$ws.Range("O12").Value = 0.0008374991573700807
$ws.Range("P12").Value = 2.119553001521041
$ws.Range("Q12").Value = "[1.5660792207084242, 2.6730267823336575]"
$ws.Range("R12").Value = [double]"9.045129090168302e-10"
$ws.Range("S12").Value = [double]"9.045129090168302e-10"
$ws.Range("T12").Value = 47.9809949382835
$ws.Range("U12").Value = "[32.968799208525056, 62.99319066804195]"
$ws.Range("V12").Value = [double]"6.987197642693843e-08"
$ws.Range("W12").Value = [double]"6.987197642693843e-08"
$ws.Range("X12").Value = 16.36114114114142
$ws.Range("Y12").Value = 14.18624624624649
$ws.Range("Z12").Value = 18.53603603603636

# Row 13
$ws.Range("F13").Value = 24.69000000000042
$ws.Range("H13").Value = [double]"7.42185693198838e-05"
$ws.Range("I13").Value = [double]"7.42185693198838e-05"
$ws.Range("L13").Value = 46.56776880799469
$ws.Range("M13").Value = "[25.625217327387844, 67.51032028860153]"
$ws.Range("N13").Value = [double]"5.1022483644525e-05"
$ws.Range("O13").Value = [double]"5.1022483644525e-05"
$ws.Range("P13").Value = 1.591237119836271
$ws.Range("Q13").Value = "[0.9748685912040393, 2.207605648468503]"
$ws.Range("R13").Value = [double]"4.718438163964223e-06"
$ws.Range("S13").Value = [double]"4.718438163964223e-06"
$ws.Range("T13").Value = 65.75273610225001
$ws.Range("U13").Value = "[51.94221420904762, 79.5632579954524]"
$ws.Range("V13").Value = [double]"1.902700219602593e-12"
$ws.Range("W13").Value = [double]"1.902700219602593e-12"
$ws.Range("X13").Value = 18.43717717717749
$ws.Range("Y13").Value = 16.01513513513541
$ws.Range("Z13").Value = 20.85921921921958

# Row 14
$ws.Range("F14").Value = 24.69000000000042
$ws.Range("H14").Value = 0.002303742582019264
$ws.Range("I14").Value = 0.002303742582019264
$ws.Range("L14").Value = 38.2934716691858
$ws.Range("M14").Value = "[12.74280790406393, 63.84413543430767]"
$ws.Range("N14").Value = 0.004172090758335578
$ws.Range("O14").Value = 0.004172090758335578
$ws.Range("P14").Value = 1.754763464167272
$ws.Range("Q14").Value = "[0.9497106920761933, 2.559816236258351]"
$ws.Range("R14").Value = [double]"6.785998617564459e-05"
$ws.Range("S14").Value = [double]"6.785998617564459e-05"
$ws.Range("T14").Value = 61.4867125652187
$ws.Range("U14").Value = "[46.67293094267262, 76.30049418776477]"
$ws.Range("V14").Value = [double]"1.038267249953151e-10"
$ws.Range("W14").Value = [double]"1.038267249953151e-10"
$ws.Range("X14").Value = 17.7945945945949
$ws.Range("Y14").Value = 14.63111111111136
$ws.Range("Z14").Value = 20.95807807807844
